$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.384.02"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").Value = "2.644.80"
$ws.Range("E3").Value = "  +0.46%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.05"
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.85"
$ws.Range("E6").Value = "  +0.66%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.547"
$ws.Range("E8").Value = "  -0.21%  "
$ws.Range("D9").Value = "2.642.82"
$ws.Range("E9").Value = "  +0.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.145"
$ws.Range("E10").Value = "  +8.13%  "
$ws.Range("E11").Value = "  -0.62%  "
$ws.Range("E12").Value = "  +1.14%  "
$ws.Range("E13").Value = "  +1.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.23"
$ws.Range("E14").Value = "  +1.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000192"
$ws.Range("E15").Value = "  +1.83%  "
$ws.Range("D16").Value = "3.122.77"
$ws.Range("E16").Value = "  +0.32%  "
$ws.Range("D17").Value = "68.331.31"
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("D18").Value = "2.641.95"
$ws.Range("E18").Value = "  +0.39%  "
$ws.Range("E19").Value = "  +1.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "364.92"
$ws.Range("E20").Value = "  -2.63%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.52"
$ws.Range("E21").Value = "  +0.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.40"
$ws.Range("E22").Value = "  +3.40%  "
$ws.Range("E23").Value = "  +1.39%  "
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.45"
$ws.Range("E25").Value = "  +2.20%  "
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.81"
$ws.Range("E27").Value = "  -1.65%  "
$ws.Range("E28").Value = "  +1.26%  "
$ws.Range("D29").Value = "2.775.40"
$ws.Range("E29").Value = "  +0.43%  "
$ws.Range("E30").Value = "  +0.35%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "573.34"
$ws.Range("E31").Value = "  -1.35%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.21"
$ws.Range("E32").Value = "  +4.48%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.42"
$ws.Range("E33").Value = "  +1.24%  "
$ws.Range("E34").Value = "  +0.45%  "
$ws.Range("E35").Value = "  +3.15%  "
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("E37").Value = "  +4.91%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "160.84"
$ws.Range("E38").Value = "  +1.19%  "
$ws.Range("E39").Value = "  +1.09%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.91"
$ws.Range("E40").Value = "  +0.32%  "
$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.374"
$ws.Range("E41").Value = "  +1.32%  "
$ws.Range("E42").Value = "  +1.55%  "
$ws.Range("D43").Value = "0.0₆0339"
$ws.Range("E43").Value = "  +6.71%  "
$ws.Range("E44").Value = "  +0.17%  "
$ws.Range("E45").Value = "  +3.55%  "
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.38"
$ws.Range("E47").Value = "  -0.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "156.99"
$ws.Range("E48").Value = "  +0.63%  "
$ws.Range("E49").Value = "  +1.91%  "
$ws.Range("E50").Value = "  +1.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.99"
$ws.Range("E51").Value = "  +0.45%  "
